$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "40.128.45"
$ws.Range("E2").Value = "  +0.20%  "
$ws.Range("D3").Value = "2.226.78"
$ws.Range("E3").Value = "  +0.61%  "
$ws.Range("D5").Value = "'294.16"
$ws.Range("E5").Value = "  +1.74%  "
$ws.Range("D6").Value = "'88.15"
$ws.Range("E6").Value = "  +0.31%  "
$ws.Range("E7").Value = "  -0.25%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("E9").Value = "  -0.17%  "
$ws.Range("D10").Value = "'30.86"
$ws.Range("E10").Value = "  +0.61%  "
$ws.Range("D11").Value = "'50.92"
$ws.Range("E11").Value = "  +6.32%  "
$ws.Range("E12").Value = "  +0.14%  "
$ws.Range("E13").Value = "  +3.09%  "
$ws.Range("E14").Value = "  -0.03%  "
$ws.Range("D15").Value = "2.586.96"
$ws.Range("D16").Value = "'13.88"
$ws.Range("E16").Value = "  -1.00%  "
$ws.Range("D17").Value = "2.253.94"
$ws.Range("E17").Value = "  +2.55%  "
$ws.Range("E18").Value = "  +1.47%  "
$ws.Range("D19").Value = "40.091.12"
$ws.Range("E19").Value = "  +0.31%  "
$ws.Range("E20").Value = "  +0.47%  "
$ws.Range("D21").Value = "'11.36"
$ws.Range("E21").Value = "  -4.86%  "
$ws.Range("E22").Value = "  -0.16%  "
$ws.Range("D23").Value = "'65.77"
$ws.Range("D24").Value = "'236.83"
$ws.Range("E24").Value = "  +0.56%  "
$ws.Range("E26").Value = "  +0.95%  "
$ws.Range("D27").Value = "'1.84"
$ws.Range("E27").Value = "  -0.38%  "
$ws.Range("D28").Value = "'23.27"
$ws.Range("E28").Value = "  +2.84%  "
$ws.Range("D29").Value = "'9.34"
$ws.Range("E29").Value = "  +1.14%  "
$ws.Range("E30").Value = "  -6.53%  "
$ws.Range("D31").Value = "'158.82"
$ws.Range("E31").Value = "  +3.95%  "
$ws.Range("D32").Value = "'31.99"
$ws.Range("E32").Value = "  -0.66%  "
$ws.Range("E33").Value = "  -0.05%  "
$ws.Range("E34").Value = "  +0.63%  "
$ws.Range("D35").Value = "'3.02"
$ws.Range("E35").Value = "  +6.66%  "
$ws.Range("D36").Value = "'0.0718"
$ws.Range("E36").Value = "  +0.02%  "
$ws.Range("D37").Value = "'2.31"
$ws.Range("E37").Value = "  -2.87%  "
$ws.Range("E38").Value = "  +1.47%  "
$ws.Range("E39").Value = "  +3.22%  "
$ws.Range("D40").Value = "'0.0996"
$ws.Range("E40").Value = "  -0.07%  "
$ws.Range("D41").Value = "'15.60"
$ws.Range("E41").Value = "  -2.40%  "
$ws.Range("D42").Value = "2.089.78"
$ws.Range("E42").Value = "  +0.14%  "
$ws.Range("D43").Value = "'3.75"
$ws.Range("E43").Value = "  -1.82%  "
$ws.Range("D44").Value = "'19.08"
$ws.Range("E44").Value = "  +8.92%  "
$ws.Range("D45").Value = "'10.11"
$ws.Range("E45").Value = "  +2.41%  "
$ws.Range("E46").Value = "  +0.70%  "
$ws.Range("E47").Value = "  +3.22%  "
$ws.Range("E48").Value = "  -11.52%  "
$ws.Range("D49").Value = "2.452.93"
$ws.Range("E49").Value = "  +1.04%  "
$ws.Range("E50").Value = "  +2.92%  "
$ws.Range("E51").Value = "  +3.87%  "
